{"js": "// Remove the guidelines about creating a server to proxy API calls.\n//\n// 1. Delete the paragraph \"To handle API requests we will need to build a\n//    simple server that would be responsible for:\" and the two bullet\n//    items underneath it (\"Creating a request: ...\" and \"Handling the\n//    response, ...\").\n// 2. The trailing \"_GoBack\" bookmark (left behind by Word at the very end\n//    of the document, after the last edit) moves up to sit right after\n//    the \"It's free RESTful API...\" paragraph, i.e. where the cursor was\n//    when the removed content was last edited.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst targetTexts = [\n  \"To handle API requests we will need to build a simple server that would be responsible for:\",\n  \"Creating a request: formatting the URL and parameters\",\n  \"Handling the response, possibly performing some data mappings\",\n];\n\nconst anchorText = \"It\\u2019s free RESTful API which provides JSON data.\";\n\nlet anchorParagraph = null;\nconst toDelete = [];\n\nfor (const paragraph of items) {\n  if (paragraph.text === anchorText) {\n    anchorParagraph = paragraph;\n  } else if (targetTexts.includes(paragraph.text)) {\n    toDelete.push(paragraph);\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\"Could not find the anchor paragraph ending in 'JSON data.'\");\n}\nif (toDelete.length !== targetTexts.length) {\n  throw new Error(\"Expected to find \" + targetTexts.length + \" paragraphs to delete, found \" + toDelete.length);\n}\n\n// Remove the proxy-server paragraphs.\nfor (const paragraph of toDelete) {\n  paragraph.delete();\n}\n\n// Drop the old \"_GoBack\" bookmark wherever it currently lives...\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// ...and re-create it at the end of the \"It's free RESTful API...\"\n// paragraph, which is where it now belongs.\nconst endOfAnchor = anchorParagraph.getRange(\"End\");\nendOfAnchor.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Remove the guidelines about creating a server to proxy API calls.\n#\n# 1. Delete the paragraph \"To handle API requests we will need to build a\n#    simple server that would be responsible for:\" and the two bullet\n#    items underneath it (\"Creating a request: ...\" and \"Handling the\n#    response, ...\").\n# 2. The trailing \"_GoBack\" bookmark (left behind by Word at the very end\n#    of the document, after the last edit) moves up to sit right after\n#    the \"It's free RESTful API...\" paragraph, i.e. where the cursor was\n#    when the removed content was last edited.\n\n$d = $word.ActiveDocument\n\n$targetTexts = @(\n    \"To handle API requests we will need to build a simple server that would be responsible for:\",\n    \"Creating a request: formatting the URL and parameters\",\n    \"Handling the response, possibly performing some data mappings\"\n)\n\n# Delete the three proxy-server paragraphs (search & remove by exact text).\n# Paragraph objects captured before a mutation can go stale once earlier\n# content shifts, so re-scan from the live collection for every delete\n# instead of caching references up front.\nforeach ($targetText in $targetTexts) {\n    $deleted = $false\n    $total = $d.Paragraphs.Count\n    for ($i = 1; $i -le $total; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($text -eq $targetText) {\n            $p.Range.Delete()\n            $deleted = $true\n            break\n        }\n    }\n    if (-not $deleted) {\n        throw \"Could not find paragraph to delete: $targetText\"\n    }\n}\n\n# Drop the old \"_GoBack\" bookmark wherever it currently lives...\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# ...and re-create it at the end of the \"It's free RESTful API...\"\n# paragraph (just before the paragraph mark), which is where it now\n# belongs. Re-locate the paragraph now (after the deletes above) using an\n# ASCII-safe substring match (the source text contains a curly\n# apostrophe), so the range reflects the document's current state.\n$anchorParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -like \"*free RESTful API which provides JSON data.*\") {\n        $anchorParagraph = $p\n        break\n    }\n}\nif ($null -eq $anchorParagraph) {\n    throw \"Could not find the anchor paragraph ending in 'JSON data.'\"\n}\n\n# A bookmark collapsed exactly onto a paragraph-mark position cannot be\n# created directly (Bookmarks.Add silently resets such a range to the\n# start of the document), so nudge it into place: temporarily append a\n# placeholder character after the paragraph's text, anchor the collapsed\n# bookmark right before that placeholder (a perfectly ordinary text\n# position), then remove the placeholder. The bookmark stays put, now\n# sitting right at the end of the paragraph's text as intended.\n$endOfText = $anchorParagraph.Range.End - 1\n$insertPoint = $d.Range($endOfText, $endOfText)\n$insertPoint.InsertAfter([char]1)\n$bookmarkSpot = $d.Range($endOfText, $endOfText)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkSpot)\n$placeholderRange = $d.Range($endOfText, $endOfText + 1)\n$placeholderRange.Delete()\n"}
